$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 597.02
$ws.Range("I15").Value = 597.02
$ws.Range("K15").Value = 1791.06
$ws.Range("M15").Value = -1622.06
$ws.Range("H112").Value = 16130566
$ws.Range("J112").Value = 1564.3334
$ws.Range("L112").Value = 4693.0002
$ws.Range("N112").Value = -6909.0002
$ws.Range("H129").Value = 746.1875
$ws.Range("I129").Value = 400.7143
$ws.Range("J129").Value = 1014.8889
$ws.Range("K129").Value = 1202.1429
$ws.Range("L129").Value = 3044.6667
$ws.Range("M129").Value = 3797.8571
$ws.Range("N129").Value = -13044.6667
$ws.Range("H132").Value = 27033062
$ws.Range("I132").Value = 33339284
$ws.Range("J132").Value = 6401.7144
$ws.Range("K132").Value = 100017852
$ws.Range("L132").Value = 19205.1432
$ws.Range("M132").Value = -100015322
$ws.Range("N132").Value = -24265.1432
$ws.Range("H137").Value = 3724.8542
$ws.Range("I137").Value = 3343.0386
$ws.Range("J137").Value = 4176.091
$ws.Range("K137").Value = 10029.1158
$ws.Range("L137").Value = 12528.273
$ws.Range("M137").Value = -7479.1158
$ws.Range("N137").Value = -17628.273
$ws.Range("H138").Value = 5096.37
$ws.Range("I138").Value = 649.7222
$ws.Range("J138").Value = 6072.4634
$ws.Range("K138").Value = 1949.1666
$ws.Range("L138").Value = 18217.3902
$ws.Range("M138").Value = 3190.8334
$ws.Range("N138").Value = -28497.3902

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1447.99
$ws.Range("I32").Value = 1123.5679
$ws.Range("J32").Value = 2831.0527
$ws.Range("K32").Value = 1123.5679
$ws.Range("L32").Value = 2831.0527
$ws.Range("M32").Value = -836.5679
$ws.Range("N32").Value = -3405.0527
$ws.Range("H122").Value = 1681.8
$ws.Range("I122").Value = 1185.8055
$ws.Range("J122").Value = 3665.7778
$ws.Range("K122").Value = 3557.4165
$ws.Range("L122").Value = 10997.3334
$ws.Range("M122").Value = -1107.4165
$ws.Range("N122").Value = -15897.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1775.3835
$ws.Range("I134").Value = 990.322
$ws.Range("J134").Value = 5083.857
$ws.Range("K134").Value = 2970.966
$ws.Range("L134").Value = 15251.571
$ws.Range("M134").Value = -435.9659999999999
$ws.Range("N134").Value = -20321.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9093526
$ws.Range("I31").Value = 1340.909
$ws.Range("J31").Value = 22731804
$ws.Range("K31").Value = 1340.909
$ws.Range("L31").Value = 22731804
$ws.Range("M31").Value = -1045.909
$ws.Range("N31").Value = -22732394
$ws.Range("H34").Value = 9093526
$ws.Range("I34").Value = 1340.909
$ws.Range("J34").Value = 22731804
$ws.Range("K34").Value = 1340.909
$ws.Range("L34").Value = 22731804
$ws.Range("M34").Value = -1138.909
$ws.Range("N34").Value = -22732208
$ws.Range("H58").Value = 1776.0919
$ws.Range("I58").Value = 1623.1177
$ws.Range("J58").Value = 2323.5789
$ws.Range("K58").Value = 1623.1177
$ws.Range("L58").Value = 2323.5789
$ws.Range("M58").Value = -1420.1177
$ws.Range("N58").Value = -2729.5789
$ws.Range("H94").Value = 1993.4546
$ws.Range("J94").Value = 1993.4546
$ws.Range("L94").Value = 1993.4546
$ws.Range("N94").Value = -2895.4546
$ws.Range("H107").Value = 686.2778
$ws.Range("I107").Value = 331.9091
$ws.Range("J107").Value = 1243.1428
$ws.Range("K107").Value = 331.9091
$ws.Range("L107").Value = 1243.1428
$ws.Range("M107").Value = 1588.0909
$ws.Range("N107").Value = -5083.1428
$ws.Range("H132").Value = 1963.2264
$ws.Range("I132").Value = 1861.8718
$ws.Range("J132").Value = 2245.5715
$ws.Range("K132").Value = 5585.6154
$ws.Range("L132").Value = 6736.7145
$ws.Range("M132").Value = -3055.6154
$ws.Range("N132").Value = -11796.7145
$ws.Range("H136").Value = 1776.0919
$ws.Range("I136").Value = 1623.1177
$ws.Range("J136").Value = 2323.5789
$ws.Range("K136").Value = 4869.3531
$ws.Range("L136").Value = 6970.736699999999
$ws.Range("M136").Value = -2319.3531
$ws.Range("N136").Value = -12070.7367

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 9100
$ws.Range("I80").Value = 7998
$ws.Range("J80").Value = 9237.75
$ws.Range("K80").Value = 23994
$ws.Range("L80").Value = 27713.25
$ws.Range("M80").Value = -23058
$ws.Range("N80").Value = -29585.25
$ws.Range("H83").Value = 9100
$ws.Range("I83").Value = 7998
$ws.Range("J83").Value = 9237.75
$ws.Range("K83").Value = 71982
$ws.Range("L83").Value = 83139.75
$ws.Range("M83").Value = -67302
$ws.Range("N83").Value = -92499.75
$ws.Range("H98").Value = 287.76923
$ws.Range("I98").Value = 92.333336
$ws.Range("J98").Value = 346.4
$ws.Range("K98").Value = 277.000008
$ws.Range("L98").Value = 1039.2
$ws.Range("M98").Value = 1220.999992
$ws.Range("N98").Value = -4035.2
$ws.Range("H113").Value = 612.4
$ws.Range("I113").Value = 513.907
$ws.Range("J113").Value = 965.3333
$ws.Range("K113").Value = 1541.721
$ws.Range("L113").Value = 2895.9999
$ws.Range("M113").Value = 628.279
$ws.Range("N113").Value = -7235.9999
$ws.Range("H131").Value = 988.6429000000001
$ws.Range("J131").Value = 1076.4681
$ws.Range("L131").Value = 3229.4043
$ws.Range("N131").Value = -13309.4043

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 49800
$ws.Range("J68").Value = 49800
$ws.Range("L68").Value = 49800
$ws.Range("N68").Value = -51422
$ws.Range("H71").Value = 49800
$ws.Range("J71").Value = 49800
$ws.Range("L71").Value = 149400
$ws.Range("N71").Value = -157512
$ws.Range("H107").Value = 578.1053000000001
$ws.Range("I107").Value = 338.9
$ws.Range("J107").Value = 843.8889
$ws.Range("K107").Value = 338.9
$ws.Range("L107").Value = 843.8889
$ws.Range("M107").Value = 1581.1
$ws.Range("N107").Value = -4683.8889
$ws.Range("H132").Value = 2395.8333
$ws.Range("I132").Value = 1638.1818
$ws.Range("J132").Value = 3586.4285
$ws.Range("K132").Value = 4914.5454
$ws.Range("L132").Value = 10759.2855
$ws.Range("M132").Value = -2384.5454
$ws.Range("N132").Value = -15819.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5742.3076
$ws.Range("I40").Value = 4003.4075
$ws.Range("J40").Value = 9654.833000000001
$ws.Range("K40").Value = 4003.4075
$ws.Range("L40").Value = 9654.833000000001
$ws.Range("M40").Value = -3867.4075
$ws.Range("N40").Value = -9926.833000000001
$ws.Range("H46").Value = 2182.2778
$ws.Range("I46").Value = 1697.6666
$ws.Range("J46").Value = 2666.889
$ws.Range("K46").Value = 1697.6666
$ws.Range("L46").Value = 2666.889
$ws.Range("M46").Value = -1509.6666
$ws.Range("N46").Value = -3042.889
$ws.Range("H68").Value = 855.1892
$ws.Range("I68").Value = 710.05884
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 710.05884
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = 38.94115999999997
$ws.Range("N68").Value = -3998
$ws.Range("H71").Value = 855.1892
$ws.Range("I71").Value = 710.05884
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 3550.2942
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = 193.7057999999997
$ws.Range("N71").Value = -19988
$ws.Range("H74").Value = 34653.223
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 37734.875
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 37734.875
$ws.Range("M74").Value = -9002
$ws.Range("N74").Value = -39730.875
$ws.Range("H77").Value = 34653.223
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 37734.875
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 113204.625
$ws.Range("M77").Value = -25008
$ws.Range("N77").Value = -123188.625
$ws.Range("H82").Value = 3820.875
$ws.Range("I82").Value = 6196.222
$ws.Range("J82").Value = 1877.409
$ws.Range("K82").Value = 6196.222
$ws.Range("L82").Value = 1877.409
$ws.Range("M82").Value = -5835.222
$ws.Range("N82").Value = -2599.409
$ws.Range("H85").Value = 3820.875
$ws.Range("I85").Value = 6196.222
$ws.Range("J85").Value = 1877.409
$ws.Range("K85").Value = 6196.222
$ws.Range("L85").Value = 1877.409
$ws.Range("M85").Value = -4948.222
$ws.Range("N85").Value = -4373.409
$ws.Range("H92").Value = 32759.334
$ws.Range("J92").Value = 32759.334
$ws.Range("L92").Value = 32759.334
$ws.Range("N92").Value = -37751.334
$ws.Range("H122").Value = 5250.5
$ws.Range("I122").Value = 1938.375
$ws.Range("J122").Value = 9666.666999999999
$ws.Range("K122").Value = 5815.125
$ws.Range("L122").Value = 29000.001
$ws.Range("M122").Value = -3365.125
$ws.Range("N122").Value = -33900.001
$ws.Range("H132").Value = 3263.6892
$ws.Range("I132").Value = 1100.5094
$ws.Range("J132").Value = 8723.143
$ws.Range("K132").Value = 3301.5282
$ws.Range("L132").Value = 26169.429
$ws.Range("M132").Value = -771.5281999999997
$ws.Range("N132").Value = -31229.429
$ws.Range("H136").Value = 2156.638
$ws.Range("I136").Value = 1325.738
$ws.Range("J136").Value = 4337.75
$ws.Range("K136").Value = 3977.214
$ws.Range("L136").Value = 13013.25
$ws.Range("M136").Value = -1427.214
$ws.Range("N136").Value = -18113.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 24124
$ws.Range("I75").Value = 8118
$ws.Range("J75").Value = 40130
$ws.Range("K75").Value = 8118
$ws.Range("L75").Value = 40130
$ws.Range("M75").Value = -7182
$ws.Range("N75").Value = -42002
$ws.Range("H78").Value = 24124
$ws.Range("I78").Value = 8118
$ws.Range("J78").Value = 40130
$ws.Range("K78").Value = 24354
$ws.Range("L78").Value = 120390
$ws.Range("M78").Value = -19674
$ws.Range("N78").Value = -129750
$ws.Range("H81").Value = 24727320
$ws.Range("I81").Value = 29222688
$ws.Range("J81").Value = 2800
$ws.Range("K81").Value = 58445376
$ws.Range("L81").Value = 5600
$ws.Range("M81").Value = -58444315
$ws.Range("N81").Value = -7722
$ws.Range("H84").Value = 24727320
$ws.Range("I84").Value = 29222688
$ws.Range("J84").Value = 2800
$ws.Range("K84").Value = 292226880
$ws.Range("L84").Value = 28000
$ws.Range("M84").Value = -292221576
$ws.Range("N84").Value = -38608
$ws.Range("H130").Value = 37391
$ws.Range("J130").Value = 37391
$ws.Range("L130").Value = 37391
$ws.Range("N130").Value = -47431
$ws.Range("H132").Value = 5130966.5
$ws.Range("I132").Value = 2574.9038
$ws.Range("J132").Value = 25644532
$ws.Range("K132").Value = 7724.7114
$ws.Range("L132").Value = 76933596
$ws.Range("M132").Value = -5194.7114
$ws.Range("N132").Value = -76938656
